$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32, shifting existing rows 32-36 down to 33-37
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new record's data
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44543
$ws.Range("D32").NumberFormat = $ws.Range("D33").NumberFormat
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112022
$ws.Range("G32").Value = "Arveja Verde"
$ws.Range("H32").Value = "Perfection"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13500
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Provincia de Diguillín"
$ws.Range("P32").Value = 540
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
